# Scenario_HeatingSystem_MinimumRenewablePercentage.xlsx
# - design of scenario-specific developments - calibration of person number for population

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the minimum-renewable-percentage scenario values ---------------
# Scenario 10 (id_scenario=10, rows 19-35): "65% renewable heating after 2024 (GEG)"
#   2025..2050 (cols U..AT) drop from 0.8 to 0.65 (2024 / col T already 0.65).
$ws.Range("U19:AT35").Value = 0.65

# Scenario 20 (id_scenario=20, rows 36-52): "65% renewable heating after 2030 (late GEG)"
#   2024..2029 (cols T..Y) drop from 0.65 to 0 ; 2030 onward (Z..AT) unchanged.
$ws.Range("T36:Y52").Value = 0

# Scenario 30 (id_scenario=30, rows 53-69): "65% renewable heating after 2035 (late GEG)"
#   2024..2034 (cols T..AD) drop from 0.65 to 0 ; 2035 onward (AE..AT) unchanged.
$ws.Range("T53:AD69").Value = 0

# --- Add the new "note" lookup sheet after Sheet1 ---------------------------
$note = $wb.Worksheets.Add($null, $ws)
$note.Name = "note"

$note.Range("A1").Value = "id_scenario"
$note.Range("B1").Value = "note"

$note.Range("A2").Value = 10
$note.Range("B2").Value = "65% renewable heating after 2024 (GEG)"

$note.Range("A3").Value = 20
$note.Range("B3").Value = "65% renewable heating after 2030 (late implementation of GEG)"

$note.Range("A4").Value = 30
$note.Range("B4").Value = "65% renewable heating after 2035 (late implementation of GEG)"

$note.Range("A5").Value = 1
$note.Range("B5").Value = "65% renewable heating after 2024 (GEG)"

# --- Update the sheet view / selection on Sheet1 and re-activate it --------
$ws.Activate()
$ws.Range("T19:AT35").Select()
